# Update NATMI TPM-derived ligand/receptor expression and specificity values
# (columns G-J: ligand stats by sending cluster; M-P: receptor stats by target
#  cluster; Q-T: derived edge weights/specificities) to match new TPM inputs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1194.705281333333
$ws.Range("H2").Value = 3584.115844
$ws.Range("I2").Value = 0.362469594586466
$ws.Range("J2").Value = 0.362469594586466
$ws.Range("M2").Value = 24.576554
$ws.Range("N2").Value = 73.729662
$ws.Range("O2").Value = 0.07553767049546639
$ws.Range("P2").Value = 0.07553767049546638
$ws.Range("Q2").Value = 29361.73886077386
$ws.Range("R2").Value = 264255.6497469647
$ws.Range("S2").Value = 0.02738010880049776
$ws.Range("T2").Value = 0.02738010880049776

# Row 3
$ws.Range("G3").Value = 1194.705281333333
$ws.Range("H3").Value = 3584.115844
$ws.Range("I3").Value = 0.362469594586466
$ws.Range("J3").Value = 0.362469594586466
$ws.Range("O3").Value = 0.359764849016532
$ws.Range("P3").Value = 0.359764849016532
$ws.Range("Q3").Value = 139841.7700575388
$ws.Range("R3").Value = 1258575.93051785
$ws.Range("S3").Value = 0.1304038189694835
$ws.Range("T3").Value = 0.1304038189694835

# Row 4
$ws.Range("G4").Value = 1194.705281333333
$ws.Range("H4").Value = 3584.115844
$ws.Range("I4").Value = 0.362469594586466
$ws.Range("J4").Value = 0.362469594586466
$ws.Range("M4").Value = 55.68784966666667
$ws.Range("N4").Value = 167.063549
$ws.Range("O4").Value = 0.1711603033819035
$ws.Range("P4").Value = 0.1711603033819035
$ws.Range("Q4").Value = 66530.56810286337
$ws.Range("R4").Value = 598775.1129257705
$ws.Range("S4").Value = 0.0620404057761351
$ws.Range("T4").Value = 0.0620404057761351

# Row 5
$ws.Range("G5").Value = 1194.705281333333
$ws.Range("H5").Value = 3584.115844
$ws.Range("I5").Value = 0.362469594586466
$ws.Range("J5").Value = 0.362469594586466
$ws.Range("M5").Value = 128.0392633333333
$ws.Range("N5").Value = 384.11779
$ws.Range("O5").Value = 0.3935371771060981
$ws.Range("P5").Value = 0.3935371771060981
$ws.Range("Q5").Value = 152969.1841223628
$ws.Range("R5").Value = 1376722.657101265
$ws.Range("S5").Value = 0.1426452610403497
$ws.Range("T5").Value = 0.1426452610403497

# Row 6
$ws.Range("I6").Value = 0.3148943969447861
$ws.Range("J6").Value = 0.314894396944786
$ws.Range("M6").Value = 24.576554
$ws.Range("N6").Value = 73.729662
$ws.Range("O6").Value = 0.07553767049546639
$ws.Range("P6").Value = 0.07553767049546638
$ws.Range("Q6").Value = 25507.92449877642
$ws.Range("R6").Value = 229571.3204889878
$ws.Range("S6").Value = 0.02378638919728385
$ws.Range("T6").Value = 0.02378638919728384

# Row 7
$ws.Range("I7").Value = 0.3148943969447861
$ws.Range("J7").Value = 0.314894396944786
$ws.Range("O7").Value = 0.359764849016532
$ws.Range("P7").Value = 0.359764849016532
$ws.Range("S7").Value = 0.1132879351729929
$ws.Range("T7").Value = 0.1132879351729928

# Row 8
$ws.Range("I8").Value = 0.3148943969447861
$ws.Range("J8").Value = 0.314894396944786
$ws.Range("M8").Value = 55.68784966666667
$ws.Range("N8").Value = 167.063549
$ws.Range("O8").Value = 0.1711603033819035
$ws.Range("P8").Value = 0.1711603033819035
$ws.Range("Q8").Value = 57798.23586319486
$ws.Range("R8").Value = 520184.1227687537
$ws.Range("S8").Value = 0.05389742051433114
$ws.Range("T8").Value = 0.05389742051433113

# Row 9
$ws.Range("I9").Value = 0.3148943969447861
$ws.Range("J9").Value = 0.314894396944786
$ws.Range("M9").Value = 128.0392633333333
$ws.Range("N9").Value = 384.11779
$ws.Range("O9").Value = 0.3935371771060981
$ws.Range("P9").Value = 0.3935371771060981
$ws.Range("Q9").Value = 132891.5299510915
$ws.Range("R9").Value = 1196023.769559824
$ws.Range("S9").Value = 0.1239226520601782
$ws.Range("T9").Value = 0.1239226520601782

# Row 10
$ws.Range("G10").Value = 560.3422443333334
$ws.Range("H10").Value = 1681.026733
$ws.Range("I10").Value = 0.170005966581565
$ws.Range("J10").Value = 0.170005966581565
$ws.Range("M10").Value = 24.576554
$ws.Range("N10").Value = 73.729662
$ws.Range("O10").Value = 0.07553767049546639
$ws.Range("P10").Value = 0.07553767049546638
$ws.Range("Q10").Value = 13771.28142633936
$ws.Range("R10").Value = 123941.5328370543
$ws.Range("S10").Value = 0.01284185468590153
$ws.Range("T10").Value = 0.01284185468590152

# Row 11
$ws.Range("G11").Value = 560.3422443333334
$ws.Range("H11").Value = 1681.026733
$ws.Range("I11").Value = 0.170005966581565
$ws.Range("J11").Value = 0.170005966581565
$ws.Range("O11").Value = 0.359764849016532
$ws.Range("P11").Value = 0.359764849016532
$ws.Range("Q11").Value = 65588.77114708623
$ws.Range("R11").Value = 590298.9403237761
$ws.Range("S11").Value = 0.06116217089912632
$ws.Range("T11").Value = 0.06116217089912632

# Row 12
$ws.Range("G12").Value = 560.3422443333334
$ws.Range("H12").Value = 1681.026733
$ws.Range("I12").Value = 0.170005966581565
$ws.Range("J12").Value = 0.170005966581565
$ws.Range("M12").Value = 55.68784966666667
$ws.Range("N12").Value = 167.063549
$ws.Range("O12").Value = 0.1711603033819035
$ws.Range("P12").Value = 0.1711603033819035
$ws.Range("Q12").Value = 31204.25466431728
$ws.Range("R12").Value = 280838.2919788555
$ws.Range("S12").Value = 0.02909827281683442
$ws.Range("T12").Value = 0.02909827281683442

# Row 13
$ws.Range("G13").Value = 560.3422443333334
$ws.Range("H13").Value = 1681.026733
$ws.Range("I13").Value = 0.170005966581565
$ws.Range("J13").Value = 0.170005966581565
$ws.Range("M13").Value = 128.0392633333333
$ws.Range("N13").Value = 384.11779
$ws.Range("O13").Value = 0.3935371771060981
$ws.Range("P13").Value = 0.3935371771060981
$ws.Range("Q13").Value = 71745.80817898668
$ws.Range("R13").Value = 645712.2736108801
$ws.Range("S13").Value = 0.06690366817970274
$ws.Range("T13").Value = 0.06690366817970274

# Row 14
$ws.Range("G14").Value = 503.070933
$ws.Range("H14").Value = 1509.212799
$ws.Range("I14").Value = 0.1526300418871828
$ws.Range("J14").Value = 0.1526300418871828
$ws.Range("M14").Value = 24.576554
$ws.Range("N14").Value = 73.729662
$ws.Range("O14").Value = 0.07553767049546639
$ws.Range("P14").Value = 0.07553767049546638
$ws.Range("Q14").Value = 12363.74995070488
$ws.Range("R14").Value = 111273.7495563439
$ws.Range("S14").Value = 0.01152931781178325
$ws.Range("T14").Value = 0.01152931781178325

# Row 15
$ws.Range("G15").Value = 503.070933
$ws.Range("H15").Value = 1509.212799
$ws.Range("I15").Value = 0.1526300418871828
$ws.Range("J15").Value = 0.1526300418871828
$ws.Range("O15").Value = 0.359764849016532
$ws.Range("P15").Value = 0.359764849016532
$ws.Range("Q15").Value = 58885.0914400446
$ws.Range("R15").Value = 529965.8229604014
$ws.Range("S15").Value = 0.05491092397492929
$ws.Range("T15").Value = 0.05491092397492929

# Row 16
$ws.Range("G16").Value = 503.070933
$ws.Range("H16").Value = 1509.212799
$ws.Range("I16").Value = 0.1526300418871828
$ws.Range("J16").Value = 0.1526300418871828
$ws.Range("M16").Value = 55.68784966666667
$ws.Range("N16").Value = 167.063549
$ws.Range("O16").Value = 0.1711603033819035
$ws.Range("P16").Value = 0.1711603033819035
$ws.Range("Q16").Value = 28014.93848857374
$ws.Range("R16").Value = 252134.4463971637
$ws.Range("S16").Value = 0.02612420427460286
$ws.Range("T16").Value = 0.02612420427460286

# Row 17
$ws.Range("G17").Value = 503.070933
$ws.Range("H17").Value = 1509.212799
$ws.Range("I17").Value = 0.1526300418871828
$ws.Range("J17").Value = 0.1526300418871828
$ws.Range("M17").Value = 128.0392633333333
$ws.Range("N17").Value = 384.11779
$ws.Range("O17").Value = 0.3935371771060981
$ws.Range("P17").Value = 0.3935371771060981
$ws.Range("Q17").Value = 64412.83166573269
$ws.Range("R17").Value = 579715.4849915942
$ws.Range("S17").Value = 0.06006559582586745
$ws.Range("T17").Value = 0.06006559582586745
